$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price/Volume columns stay as plain text so values like "61.655.92"
# are not re-interpreted as numbers/dates by Excel.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "61.655.92"
$ws.Range("E2").Value = "  +1.28%  "
$ws.Range("D3").Value = "3.395.46"
$ws.Range("E3").Value = "  +0.24%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").Value = "576.86"
$ws.Range("E5").Value = "  +0.96%  "
$ws.Range("D6").Value = "141.47"
$ws.Range("E6").Value = "  -0.16%  "
$ws.Range("E7").Value = "  +0.04%  "
$ws.Range("E8").Value = "  -0.14%  "
$ws.Range("D9").Value = "7.75"
$ws.Range("E9").Value = "  +2.86%  "
$ws.Range("E10").Value = "  -0.38%  "
$ws.Range("D11").Value = "0.388"
$ws.Range("E11").Value = "  -1.53%  "
$ws.Range("D12").Value = "3.975.28"
$ws.Range("E12").Value = "  +0.17%  "

# Row 13 and row 14 swap content: TRON moves to row 14, Avalanche moves to row 13
$ws.Range("B13").Value = "Avalanche"
$ws.Range("C13").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D13").Value = "28.48"
$ws.Range("E13").Value = "  +1.01%  "
$ws.Range("B14").Value = "TRON"
$ws.Range("C14").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D14").Value = "0.125"
$ws.Range("E14").Value = "  +0.18%  "

$ws.Range("D15").Value = "3.382.48"
$ws.Range("E15").Value = "  -0.32%  "
$ws.Range("E16").Value = "  +0.30%  "
$ws.Range("D17").Value = "61.644.43"
$ws.Range("E17").Value = "  +1.12%  "
$ws.Range("E18").Value = "  +0.16%  "
$ws.Range("D19").Value = "13.71"
$ws.Range("E19").Value = "  -0.97%  "
$ws.Range("D20").Value = "8.99"
$ws.Range("E20").Value = "  +0.27%  "
$ws.Range("D21").Value = "391.44"
$ws.Range("E21").Value = "  +1.97%  "
$ws.Range("D22").Value = "75.52"
$ws.Range("E22").Value = "  +1.48%  "
$ws.Range("D23").Value = "0.554"
$ws.Range("E23").Value = "  -0.56%  "
$ws.Range("E24").Value = "  +0.04%  "
$ws.Range("E25").Value = "  -3.79%  "
$ws.Range("D26").Value = "0.195"
$ws.Range("E26").Value = "  +9.41%  "
$ws.Range("E27").Value = "  +0.03%  "
$ws.Range("D28").Value = "7.29"
$ws.Range("E28").Value = "  -1.09%  "
$ws.Range("D29").Value = "8.04"
$ws.Range("E29").Value = "  +0.92%  "
$ws.Range("E30").Value = "  +0.46%  "
$ws.Range("E31").Value = "  -0.01%  "
$ws.Range("E32").Value = "  -3.75%  "
$ws.Range("D33").Value = "23.42"
$ws.Range("E33").Value = "  -0.37%  "
$ws.Range("D34").Value = "6.95"
$ws.Range("E34").Value = "  -0.29%  "
$ws.Range("D35").Value = "168.04"
$ws.Range("E35").Value = "  +0.29%  "
$ws.Range("D36").Value = "5.09"
$ws.Range("E36").Value = "  +2.33%  "
$ws.Range("D37").Value = "3.431.12"
$ws.Range("E37").Value = "  +0.38%  "
$ws.Range("E38").Value = "  -0.58%  "
$ws.Range("E39").Value = "  -0.23%  "
$ws.Range("D40").Value = "26.01"
$ws.Range("E40").Value = "  -5.54%  "
$ws.Range("D41").Value = "0.781"
$ws.Range("E42").Value = "  +0.19%  "
$ws.Range("E43").Value = "  -0.18%  "
$ws.Range("E44").Value = "  +2.05%  "
$ws.Range("D45").Value = "2.474.38"
$ws.Range("E45").Value = "  -0.11%  "
$ws.Range("D46").Value = "23.17"
$ws.Range("E46").Value = "  +0.77%  "
$ws.Range("D47").Value = "6.67"
$ws.Range("E47").Value = "  -2.03%  "
$ws.Range("E49").Value = "  -1.55%  "
$ws.Range("E50").Value = "  -0.44%  "
$ws.Range("E51").Value = "  -1.02%  "
